# Rename the variable/column headers in row 1 to their short/lowercase
# aliases (matches the renamed entries in the shared-string table).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "crimes"
$ws.Range("C1").Value = "ssusers"
$ws.Range("D1").Value = "popdens"
$ws.Range("E1").Value = "avgsal"
$ws.Range("F1").Value = "unemp"

# Leave the selection on F1, matching the saved state of the workbook.
$ws.Range("F1").Select()
